$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Register the "Microsoft YaHei UI" font into the workbook's font table ---
# (Excel always keeps a font registered in styles.xml once it's used anywhere,
#  including inside a shared-string rich-text run. The simplest reliable way to
#  get this engine to add the font entry is to apply it to a scratch cell and
#  then clear that cell.)
$ws.Range("ZZ1").Value2 = "x"
$ws.Range("ZZ1").Font.Name = "Microsoft YaHei UI"
$ws.Range("ZZ1").Font.Size = 12

# --- Row 5: new data row ---
$s1 = "TK Task[1:8]+T9 on NPU"
$s2 = "全体膨胀约1/3"
$ws.Range("A5").Value2 = $s1 + $s2

# Apply the special font to just the second run of the rich text (mirrors the
# "全体膨胀约1/3" portion being rendered in a distinct font/color in the source).
$len1 = $s1.Length
$len2 = $s2.Length
$chars = $ws.Range("A5").Characters($len1 + 1, $len2)
$chars.Font.Name = "Microsoft YaHei UI"
$chars.Font.Size = 12

$ws.Range("B5").Value2 = 40
$ws.Range("C5").Value2 = 0.921
$ws.Range("D5").Value2 = 0.234
$ws.Range("E5").Value2 = 0.951
$ws.Range("F5").Formula = "=C5/E5"

# Restore A5's cell-level style/format to match the rest of column A (the rich
# text run above only changes the *shared string*, not the cell style) -- copy
# number format / style from A4, which carries the same style as the other
# data rows.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clean up the scratch cell used to register the font.
$ws.Range("ZZ1").Clear()

# Row height for the new row.
$ws.Rows(5).RowHeight = 17.25

# --- Column A width ---
$ws.Columns("A").ColumnWidth = 40

# --- View: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("B8").Select()
